# Insert a new row at position 54 (pushes old rows 54..178 down to 55..179)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("54:54").Insert()

# Populate the newly inserted row 54 with its data.
$ws.Cells.Item(54, 1).Value = 7
$ws.Cells.Item(54, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(54, 3).Value = "Ñuble"
$ws.Cells.Item(54, 4).Value = 45152
$ws.Cells.Item(54, 5).Value = 16
$ws.Cells.Item(54, 6).Value = "Fruta"
$ws.Cells.Item(54, 7).Value = 100108
$ws.Cells.Item(54, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(54, 9).Value = 100108002
$ws.Cells.Item(54, 10).Value = "Mango"
$ws.Cells.Item(54, 11).Value = "Sin especificar"
$ws.Cells.Item(54, 12).Value = "Primera"
$ws.Cells.Item(54, 13).Value = 80
$ws.Cells.Item(54, 14).Value = 8000
$ws.Cells.Item(54, 15).Value = 8000
$ws.Cells.Item(54, 16).Value = 8000
$ws.Cells.Item(54, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(54, 18).Value = "Brasil"
$ws.Cells.Item(54, 19).Value = 2000
$ws.Cells.Item(54, 20).Value = 4
